$d = $word.ActiveDocument

# Helper behaviour note: always re-fetch a fresh end-of-document Range right before
# each insertion. Reusing the same Range object across an InsertParagraphAfter() call
# does not reliably advance its Start/End past the freshly inserted paragraph mark in
# this runtime, so grab $d.Range() again each time instead.

function Append-Paragraph([string]$text) {
    $r = $d.Range()
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $r2 = $d.Range()
    $r2.Collapse(0)
    if ($text -ne "") {
        $r2.InsertAfter($text)
    }
}

# New entries for "29/01 dinsdag:" added after the existing "!!NIET VERGETEN TE
# PUSHEN!!" paragraph at the end of the document.
Append-Paragraph ""
Append-Paragraph "29/01 dinsdag:"
Append-Paragraph "Template = jeroen"
Append-Paragraph "DB in tabellen steken = Dieter"
Append-Paragraph "Tanguy = verslag"
Append-Paragraph "Ruben en seb = github"

# Move the "_GoBack" bookmark from its old location (end of the "4. speciale link..."
# paragraph) to the end of the newly-added last paragraph ("...= github").
$endR = $d.Range()
$endR.Collapse(0)
$bmRange = $d.Range($endR.Start, $endR.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
